# Auto-generated script to update cryptos list values (Wed Aug 28 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '58.702.79'
$ws.Range('E2').Value = '  -5.15%  '

# Row 3
$ws.Range('D3').Value = '2.463.58'
$ws.Range('E3').Value = '  -4.71%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '534.10'
$ws.Range('E5').Value = '  -3.68%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.56'
$ws.Range('E6').Value = '  -6.84%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.06%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.567'
$ws.Range('E8').Value = '  -4.58%  '

# Row 9
$ws.Range('D9').Value = '2.482.95'
$ws.Range('E9').Value = '  -4.24%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0996'
$ws.Range('E10').Value = '  -4.55%  '

# Row 11
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.158'
$ws.Range('E11').Value = '  -2.02%  '

# Row 12
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.60'
$ws.Range('E12').Value = '  +2.28%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.352'
$ws.Range('E13').Value = '  -3.28%  '

# Row 14
$ws.Range('D14').Value = '2.905.87'
$ws.Range('E14').Value = '  -4.37%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.78'
$ws.Range('E15').Value = '  -6.67%  '

# Row 16
$ws.Range('D16').Value = '58.608.12'
$ws.Range('E16').Value = '  -5.16%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000137'
$ws.Range('E17').Value = '  -4.47%  '

# Row 18
$ws.Range('D18').Value = '2.486.63'
$ws.Range('E18').Value = '  -3.69%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.30'
$ws.Range('E19').Value = '  -2.59%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.30'
$ws.Range('E20').Value = '  -5.07%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '322.77'
$ws.Range('E21').Value = '  -4.78%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.995'
$ws.Range('E22').Value = '  -0.37%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.73'
$ws.Range('E23').Value = '  -5.12%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.46'
$ws.Range('E24').Value = '  -3.75%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.438'
$ws.Range('E25').Value = '  -11.49%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.994'
$ws.Range('E26').Value = '  -0.48%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.161'
$ws.Range('E27').Value = '  -4.42%  '

# Row 28
$ws.Range('D28').Value = '2.591.97'
$ws.Range('E28').Value = '  -4.18%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.72'
$ws.Range('E29').Value = '  -4.24%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.99'
$ws.Range('E30').Value = '  -0.95%  '

# Row 31
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0770'
$ws.Range('E31').Value = '  -8.34%  '

# Row 32
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.25'
$ws.Range('E32').Value = '  -4.10%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.79'
$ws.Range('E33').Value = '  -6.72%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.04%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.16'
$ws.Range('E35').Value = '  -1.42%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.40'
$ws.Range('E36').Value = '  -1.18%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.48'
$ws.Range('E37').Value = '  -4.05%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.40'
$ws.Range('E38').Value = '  -6.40%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.60'
$ws.Range('E39').Value = '  -10.77%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.75'
$ws.Range('E40').Value = '  -5.16%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '306.57'
$ws.Range('E41').Value = '  -9.35%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.43'
$ws.Range('E42').Value = '  -2.79%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.69'
$ws.Range('E43').Value = '  -5.69%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.807'
$ws.Range('E44').Value = '  -9.90%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.996'
$ws.Range('E45').Value = '  -0.18%  '

# Row 46
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.78'
$ws.Range('E46').Value = '  -1.39%  '

# Row 47
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.591'
$ws.Range('E47').Value = '  -2.76%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '123.58'
$ws.Range('E48').Value = '  -1.02%  '

# Row 49
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0926'
$ws.Range('E49').Value = '  -4.08%  '

# Row 50
$ws.Range('B50').Value = 'Hedera'
$ws.Range('C50').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0520'
$ws.Range('E50').Value = '  -4.83%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0227'
$ws.Range('E51').Value = '  -5.48%  '
